$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instructions")

$ws1.Range("A3").Value = "v2 - 04.03.2021"
$ws1.Range("A5").Value = "General"
$ws1.Range("B6").Value = "All subscriptions are newly created and are therefore completely fresh"
$ws1.Range("B7").Value = "All subscriptions are given a unique name (following are samples, but always ensure that in the two tabs both team name columns are matching)"
$ws1.Range("B10").Value = "All Azure accounts"
$ws1.Range("C11").Value = "get unique names"
$ws1.Range("B17").Value = "Creation of 12 Team Azure subscriptions"
$ws1.Range("C18").Value = "Setting the spending / notification limit: `$ 100"
$ws1.Range("B20").Value = "Creation of a Service Principal per team, assignment as Subscription Owner"
$ws1.Range("B19").Value = "Creation of 5 Azure accounts per team, assignment as Subscription Contributor"
$ws1.Range("B24").Value = "Creation of 1 Admin Azure Subscription"
$ws1.Range("C25").Value = "Setting the spending / notification limit: `$ 200"
$ws1.Range("B26").Value = "Creation of 6 Azure acounts for the coaches"
$ws1.Range("C27").Value = "Assignment as Subscription Owner for Admin Subscription and all Team Subscriptions"
$ws1.Range("C28").Value = "Assignment as Admin in all Azure AD tenants (if custom AD tenants were created)"
$ws1.Range("C29").Value = "are already activated on all tenants / guest invitation is accepted"
$ws1.Range("B30").Value = "Creation of a Service Principal for the Admin Subscription (Usage of Multi-Tenant App feature)"
$ws1.Range("C31").Value = "Assignment as Subscription Owner for Admin Subscription and all Team Subscriptions"
$ws1.Range("A34").Value = "Transfer of data"
$ws1.Range("B35").Value = "Use this Excel as a template and enter the data in the `"Azure Subscriptions`" and `"User Accounts & SPs`" tabs"

# Switch the active/selected tab from "User Accounts & SPs" back to "Instructions",
# and update its selection to M17 (matches the saved view state in the target file).
$ws1.Activate() | Out-Null
$ws1.Range("M17").Select() | Out-Null

